$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Header: "Challenge 9" -> "Challenge 8"
# ------------------------------------------------------------------
$hdr = $d.Sections(1).Headers(1)
$hdrRange = $hdr.Range
$hdrRange.Find.ClearFormatting()
$hdrRange.Find.Execute("Challenge 9", $false, $false, $false, $false, $false, $true, 1, $false, "Challenge 8", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from its own (now-empty) paragraph
#    up onto the end of the "LEARNING OBJECTIVES" paragraph.
# ------------------------------------------------------------------

# Remove the existing _GoBack bookmark (it currently lives alone in its
# own paragraph right after the MAP-sensor procedure paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create it at the end of the text in the first paragraph
# ("LEARNING OBJECTIVES"). A zero-length range sitting exactly at a
# paragraph boundary cannot be handed to Bookmarks.Add directly, so we
# briefly insert a placeholder character, wrap the bookmark around it,
# then delete the placeholder -- leaving a clean, empty bookmark behind.
$objectivesPara = $d.Paragraphs(1)
$endPos = $objectivesPara.Range.End - 1
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("@")
$wrap = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $wrap)
$wrap2 = $d.Range($endPos, $endPos + 1)
$wrap2.Text = ""

# ------------------------------------------------------------------
# 3. Re-order the runs in the MAP-sensor procedure paragraph so the
#    stray " " that used to sit between "...sensor!" and "Begin with
#    the system..." is folded into the following sentence instead.
# ------------------------------------------------------------------
$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
$searchRange.Find.Execute("sensor! Begin with the system depressurized", $false, $false, $false, $false, $false, $true, 1, $false, "sensor!#Begin with the system depressurized", 2) | Out-Null

$searchRange2 = $d.Content
$searchRange2.Find.ClearFormatting()
$searchRange2.Find.Execute("#", $false, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

# ------------------------------------------------------------------
# 4. The paragraph that used to hold the bookmark is now a plain
#    empty paragraph (no bookmark, no run).
# ------------------------------------------------------------------
